$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in timing results (in nanoseconds, "4 second size") for the sorting
# methods that are finished. Merge Sort and Quick Sort are still in
# progress per the commit message, but Quick Sort's row already received a
# value here; Merge Sort (row 10) remains blank.
$ws.Range("B5").Value = 5300000    # STL::sort
$ws.Range("B7").Value = 15000      # Quick Sort
$ws.Range("B8").Value = 18100      # Insertion Sort
$ws.Range("B9").Value = 21500      # Selection Sort

# Reflect the last-edited/selected cell as in the authored workbook.
$ws.Range("B9").Select()
